$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the fluid_mass value (B3) which drives the Q_toBoil calculation (B7)
$ws.Range("B3").Value = 1.767144375
